$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.810.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.057.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.055.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.479"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.562.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.875.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.064.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "442.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.046.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("E41").Value = "  -6.09%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("E46").Value = "  +6.32%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.61%  "
